$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 104 <- old Row 105 data (match swap)
$ws.Range("B104").Value = 7331149
$ws.Range("F104").Value = 'Vaca Diez'
$ws.Range("G104").Value = 'Guabira'
$ws.Range("H104").Value = 3
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 'H'
$ws.Range("K104").Value = 1.8
$ws.Range("L104").Value = 3.6
$ws.Range("M104").Value = 3.6
$ws.Range("N104").Value = 2.3
$ws.Range("O104").Value = 3.75
$ws.Range("P104").Value = 2.8
$ws.Range("Q104").Value = 0
$ws.Range("R104").Value = 1.8
$ws.Range("S104").Value = 2
$ws.Range("T104").Value = 2.75
$ws.Range("U104").Value = 1.95
$ws.Range("V104").Value = 1.85
$ws.Range("W104").Value = 1.3
$ws.Range("X104").Value = -1
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = 0.8
$ws.Range("AA104").Value = -1
$ws.Range("AB104").Value = 0.475
$ws.Range("AC104").Value = -0.5

# Row 105 <- old Row 104 data (match swap)
$ws.Range("B105").Value = 7331148
$ws.Range("F105").Value = 'Always Ready'
$ws.Range("G105").Value = 'Real Santa Cruz'
$ws.Range("H105").Value = 2
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 'H'
$ws.Range("K105").Value = 1.2
$ws.Range("L105").Value = 6
$ws.Range("M105").Value = 9
$ws.Range("N105").Value = 1.166
$ws.Range("O105").Value = 7.5
$ws.Range("P105").Value = 16
$ws.Range("Q105").Value = -2
$ws.Range("R105").Value = 1.8
$ws.Range("S105").Value = 2
$ws.Range("T105").Value = 3.5
$ws.Range("U105").Value = 1.975
$ws.Range("V105").Value = 1.825
$ws.Range("W105").Value = 0.1659999999999999
$ws.Range("X105").Value = -1
$ws.Range("Y105").Value = -1
$ws.Range("Z105").Value = 0
$ws.Range("AA105").Value = -0
$ws.Range("AB105").Value = -1
$ws.Range("AC105").Value = 0.825

# Row 107 <- old Row 108 data (match swap)
$ws.Range("B107").Value = 7329845
$ws.Range("F107").Value = 'Universitario De Vinto'
$ws.Range("G107").Value = 'Bolivar'
$ws.Range("H107").Value = 2
$ws.Range("I107").Value = 2
$ws.Range("J107").Value = 'D'
$ws.Range("K107").Value = 3.1
$ws.Range("L107").Value = 3.4
$ws.Range("M107").Value = 2.05
$ws.Range("N107").Value = 3.4
$ws.Range("O107").Value = 3.5
$ws.Range("P107").Value = 2.05
$ws.Range("Q107").Value = 0.25
$ws.Range("R107").Value = 2.025
$ws.Range("S107").Value = 1.775
$ws.Range("T107").Value = 2.5
$ws.Range("U107").Value = 1.8
$ws.Range("V107").Value = 2
$ws.Range("W107").Value = -1
$ws.Range("X107").Value = 2.5
$ws.Range("Y107").Value = -1
$ws.Range("Z107").Value = 0.5125
$ws.Range("AA107").Value = -0.5
$ws.Range("AB107").Value = 0.8
$ws.Range("AC107").Value = -1

# Row 108 <- old Row 107 data (match swap)
$ws.Range("B108").Value = 7331150
$ws.Range("F108").Value = 'Libertad Gran Mamore FC'
$ws.Range("G108").Value = 'Royal Pari FC'
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 2
$ws.Range("J108").Value = 'A'
$ws.Range("K108").Value = 2.15
$ws.Range("L108").Value = 3.5
$ws.Range("M108").Value = 2.875
$ws.Range("N108").Value = 2.4
$ws.Range("O108").Value = 3.6
$ws.Range("P108").Value = 2.7
$ws.Range("Q108").Value = -0.25
$ws.Range("R108").Value = 2.075
$ws.Range("S108").Value = 1.725
$ws.Range("T108").Value = 2.75
$ws.Range("U108").Value = 1.95
$ws.Range("V108").Value = 1.85
$ws.Range("W108").Value = -1
$ws.Range("X108").Value = -1
$ws.Range("Y108").Value = 1.7
$ws.Range("Z108").Value = -1
$ws.Range("AA108").Value = 0.7250000000000001
$ws.Range("AB108").Value = -1
$ws.Range("AC108").Value = 0.8500000000000001

# Row 142 <- old Row 145 data (match swap)
$ws.Range("B142").Value = 7532414
$ws.Range("F142").Value = 'Independiente Petrolero'
$ws.Range("G142").Value = 'Real Santa Cruz'
$ws.Range("H142").Value = 1
$ws.Range("I142").Value = 0
$ws.Range("J142").Value = 'H'
$ws.Range("K142").Value = 1.571
$ws.Range("L142").Value = 3.75
$ws.Range("M142").Value = 5
$ws.Range("N142").Value = 1.3
$ws.Range("O142").Value = 5
$ws.Range("P142").Value = 11
$ws.Range("Q142").Value = -1.75
$ws.Range("R142").Value = 2
$ws.Range("S142").Value = 1.8
$ws.Range("T142").Value = 3
$ws.Range("U142").Value = 1.85
$ws.Range("V142").Value = 1.95
$ws.Range("W142").Value = 0.3
$ws.Range("X142").Value = -1
$ws.Range("Y142").Value = -1
$ws.Range("Z142").Value = -1
$ws.Range("AA142").Value = 0.8
$ws.Range("AB142").Value = -1
$ws.Range("AC142").Value = 0.95

# Row 143 <- old Row 144 data (match swap)
$ws.Range("B143").Value = 7532413
$ws.Range("F143").Value = 'Libertad Gran Mamore FC'
$ws.Range("G143").Value = 'Club Aurora'
$ws.Range("H143").Value = 0
$ws.Range("I143").Value = 1
$ws.Range("J143").Value = 'A'
$ws.Range("K143").Value = 2.25
$ws.Range("L143").Value = 3.3
$ws.Range("M143").Value = 2.8
$ws.Range("N143").Value = 2.375
$ws.Range("O143").Value = 3.4
$ws.Range("P143").Value = 2.875
$ws.Range("Q143").Value = -0.25
$ws.Range("R143").Value = 2.025
$ws.Range("S143").Value = 1.775
$ws.Range("T143").Value = 2.5
$ws.Range("U143").Value = 1.9
$ws.Range("V143").Value = 1.9
$ws.Range("W143").Value = -1
$ws.Range("X143").Value = -1
$ws.Range("Y143").Value = 1.875
$ws.Range("Z143").Value = -1
$ws.Range("AA143").Value = 0.7749999999999999
$ws.Range("AB143").Value = -1
$ws.Range("AC143").Value = 0.8999999999999999

# Row 144 <- old Row 143 data (match swap)
$ws.Range("B144").Value = 7532412
$ws.Range("F144").Value = 'Vaca Diez'
$ws.Range("G144").Value = 'Blooming'
$ws.Range("H144").Value = 0
$ws.Range("I144").Value = 3
$ws.Range("J144").Value = 'A'
$ws.Range("K144").Value = 1.727
$ws.Range("L144").Value = 3.75
$ws.Range("M144").Value = 4
$ws.Range("N144").Value = 2.3
$ws.Range("O144").Value = 3.6
$ws.Range("P144").Value = 2.875
$ws.Range("Q144").Value = -0.25
$ws.Range("R144").Value = 1.95
$ws.Range("S144").Value = 1.85
$ws.Range("T144").Value = 2.75
$ws.Range("U144").Value = 1.925
$ws.Range("V144").Value = 1.875
$ws.Range("W144").Value = -1
$ws.Range("X144").Value = -1
$ws.Range("Y144").Value = 1.875
$ws.Range("Z144").Value = -1
$ws.Range("AA144").Value = 0.8500000000000001
$ws.Range("AB144").Value = 0.4625
$ws.Range("AC144").Value = -0.5

# Row 145 <- old Row 142 data (match swap)
$ws.Range("B145").Value = 7532430
$ws.Range("F145").Value = 'Always Ready'
$ws.Range("G145").Value = 'Oriente Petrolero'
$ws.Range("H145").Value = 4
$ws.Range("I145").Value = 1
$ws.Range("J145").Value = 'H'
$ws.Range("K145").Value = 1.4
$ws.Range("L145").Value = 4.2
$ws.Range("M145").Value = 7
$ws.Range("N145").Value = 1.363
$ws.Range("O145").Value = 4.5
$ws.Range("P145").Value = 8.5
$ws.Range("Q145").Value = -1.5
$ws.Range("R145").Value = 2
$ws.Range("S145").Value = 1.8
$ws.Range("T145").Value = 3
$ws.Range("U145").Value = 1.9
$ws.Range("V145").Value = 1.9
$ws.Range("W145").Value = 0.363
$ws.Range("X145").Value = -1
$ws.Range("Y145").Value = -1
$ws.Range("Z145").Value = 1
$ws.Range("AA145").Value = -1
$ws.Range("AB145").Value = 0.8999999999999999
$ws.Range("AC145").Value = -1

# Row 148 <- old Row 149 data (match swap)
$ws.Range("B148").Value = 7532419
$ws.Range("F148").Value = 'Oriente Petrolero'
$ws.Range("G148").Value = 'Jorge Wilstermann'
$ws.Range("H148").Value = 3
$ws.Range("I148").Value = 0
$ws.Range("J148").Value = 'H'
$ws.Range("K148").Value = 2.2
$ws.Range("L148").Value = 2.5
$ws.Range("M148").Value = 4.5
$ws.Range("N148").Value = 2.375
$ws.Range("O148").Value = 2.45
$ws.Range("P148").Value = 4.5
$ws.Range("Q148").Value = -0.25
$ws.Range("R148").Value = 1.9
$ws.Range("S148").Value = 1.9
$ws.Range("T148").Value = 2
$ws.Range("U148").Value = 1.95
$ws.Range("V148").Value = 1.85
$ws.Range("W148").Value = 1.375
$ws.Range("X148").Value = -1
$ws.Range("Y148").Value = -1
$ws.Range("Z148").Value = 0.8999999999999999
$ws.Range("AA148").Value = -1
$ws.Range("AB148").Value = 0.95
$ws.Range("AC148").Value = -1

# Row 149 <- old Row 150 data (match swap)
$ws.Range("B149").Value = 7532420
$ws.Range("F149").Value = 'Club Aurora'
$ws.Range("G149").Value = 'Vaca Diez'
$ws.Range("H149").Value = 3
$ws.Range("I149").Value = 0
$ws.Range("J149").Value = 'H'
$ws.Range("K149").Value = 1.333
$ws.Range("L149").Value = 5
$ws.Range("M149").Value = 8
$ws.Range("N149").Value = 1.3
$ws.Range("O149").Value = 6.5
$ws.Range("P149").Value = 7
$ws.Range("Q149").Value = -1.5
$ws.Range("R149").Value = 1.8
$ws.Range("S149").Value = 2
$ws.Range("T149").Value = 3.25
$ws.Range("U149").Value = 1.95
$ws.Range("V149").Value = 1.85
$ws.Range("W149").Value = 0.3
$ws.Range("X149").Value = -1
$ws.Range("Y149").Value = -1
$ws.Range("Z149").Value = 0.8
$ws.Range("AA149").Value = -1
$ws.Range("AB149").Value = -0.5
$ws.Range("AC149").Value = 0.425

# Row 150 <- old Row 148 data (match swap)
$ws.Range("B150").Value = 7532421
$ws.Range("F150").Value = 'Guabira'
$ws.Range("G150").Value = 'Independiente Petrolero'
$ws.Range("H150").Value = 2
$ws.Range("I150").Value = 0
$ws.Range("J150").Value = 'H'
$ws.Range("K150").Value = 1.4
$ws.Range("L150").Value = 4.5
$ws.Range("M150").Value = 7.5
$ws.Range("N150").Value = 1.333
$ws.Range("O150").Value = 5.5
$ws.Range("P150").Value = 9.5
$ws.Range("Q150").Value = -1.5
$ws.Range("R150").Value = 1.85
$ws.Range("S150").Value = 1.95
$ws.Range("T150").Value = 3
$ws.Range("U150").Value = 1.825
$ws.Range("V150").Value = 1.975
$ws.Range("W150").Value = 0.333
$ws.Range("X150").Value = -1
$ws.Range("Y150").Value = -1
$ws.Range("Z150").Value = 0.8500000000000001
$ws.Range("AA150").Value = -1
$ws.Range("AB150").Value = -1
$ws.Range("AC150").Value = 0.9750000000000001

# Row 211 <- old Row 212 data (match swap)
$ws.Range("B211").Value = 8039389
$ws.Range("F211").Value = 'Royal Pari FC'
$ws.Range("G211").Value = 'Club Aurora'
$ws.Range("H211").Value = 0
$ws.Range("I211").Value = 0
$ws.Range("J211").Value = 'D'
$ws.Range("K211").Value = 2.2
$ws.Range("L211").Value = 3.3
$ws.Range("M211").Value = 2.875
$ws.Range("N211").Value = 2.625
$ws.Range("O211").Value = 3.4
$ws.Range("P211").Value = 2.625
$ws.Range("Q211").Value = 0
$ws.Range("R211").Value = 1.875
$ws.Range("S211").Value = 1.925
$ws.Range("T211").Value = 2.5
$ws.Range("U211").Value = 1.95
$ws.Range("V211").Value = 1.85
$ws.Range("W211").Value = -1
$ws.Range("X211").Value = 2.4
$ws.Range("Y211").Value = -1
$ws.Range("Z211").Value = 0
$ws.Range("AA211").Value = -0
$ws.Range("AB211").Value = -1
$ws.Range("AC211").Value = 0.8500000000000001

# Row 212 <- old Row 211 data (match swap)
$ws.Range("B212").Value = 8039390
$ws.Range("F212").Value = 'Blooming'
$ws.Range("G212").Value = 'Nacional Potosi'
$ws.Range("H212").Value = 2
$ws.Range("I212").Value = 1
$ws.Range("J212").Value = 'H'
$ws.Range("K212").Value = 2.05
$ws.Range("L212").Value = 3.6
$ws.Range("M212").Value = 3
$ws.Range("N212").Value = 1.8
$ws.Range("O212").Value = 3.75
$ws.Range("P212").Value = 4.2
$ws.Range("Q212").Value = -0.5
$ws.Range("R212").Value = 1.825
$ws.Range("S212").Value = 1.975
$ws.Range("T212").Value = 3
$ws.Range("U212").Value = 1.975
$ws.Range("V212").Value = 1.825
$ws.Range("W212").Value = 0.8
$ws.Range("X212").Value = -1
$ws.Range("Y212").Value = -1
$ws.Range("Z212").Value = 0.825
$ws.Range("AA212").Value = -1
$ws.Range("AB212").Value = 0
$ws.Range("AC212").Value = -0

# Row 214 <- old Row 215 data (match swap)
$ws.Range("B214").Value = 8038943
$ws.Range("F214").Value = 'San Jose de Oruro'
$ws.Range("G214").Value = 'Bolivar'
$ws.Range("H214").Value = 2
$ws.Range("I214").Value = 1
$ws.Range("J214").Value = 'H'
$ws.Range("K214").Value = 2.3
$ws.Range("L214").Value = 3.5
$ws.Range("M214").Value = 2.625
$ws.Range("N214").Value = 2.8
$ws.Range("O214").Value = 3.6
$ws.Range("P214").Value = 2.375
$ws.Range("Q214").Value = 0.25
$ws.Range("R214").Value = 1.8
$ws.Range("S214").Value = 2
$ws.Range("T214").Value = 3.25
$ws.Range("U214").Value = 1.975
$ws.Range("V214").Value = 1.825
$ws.Range("W214").Value = 1.8
$ws.Range("X214").Value = -1
$ws.Range("Y214").Value = -1
$ws.Range("Z214").Value = 0.8
$ws.Range("AA214").Value = -1
$ws.Range("AB214").Value = -0.5
$ws.Range("AC214").Value = 0.4125

# Row 215 <- old Row 214 data (match swap)
$ws.Range("B215").Value = 8039392
$ws.Range("F215").Value = 'Oriente Petrolero'
$ws.Range("G215").Value = 'Jorge Wilstermann'
$ws.Range("H215").Value = 2
$ws.Range("I215").Value = 1
$ws.Range("J215").Value = 'H'
$ws.Range("K215").Value = 2
$ws.Range("L215").Value = 3.25
$ws.Range("M215").Value = 3.4
$ws.Range("N215").Value = 1.727
$ws.Range("O215").Value = 4
$ws.Range("P215").Value = 4.5
$ws.Range("Q215").Value = -0.75
$ws.Range("R215").Value = 1.9
$ws.Range("S215").Value = 1.9
$ws.Range("T215").Value = 2.75
$ws.Range("U215").Value = 1.9
$ws.Range("V215").Value = 1.9
$ws.Range("W215").Value = 0.7270000000000001
$ws.Range("X215").Value = -1
$ws.Range("Y215").Value = -1
$ws.Range("Z215").Value = 0.45
$ws.Range("AA215").Value = -0.5
$ws.Range("AB215").Value = 0.45
$ws.Range("AC215").Value = -0.5
